# Adds the new "06. Phone book" solution (rows 117-120) to Sheet1,
# matching the target shared-strings + worksheet content from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell values -----------------------------------------------------------
$ws.Range("A117").Value = 'const content = Object.values(data).map(entry => `${entry.author}: ${entry.content}`).join(''\n'');'
$ws.Range("B117").Value = 'Обхождане на обект с обекти и създаване на масив със съдържанието на обектите'

$ws.Range("A118").Value = ' const content = [];
 for (const elem of Object.values(data)) {
            content.push(`${elem.author}: ${elem.content}`)
        }'
# B118 stays empty (merged into the B117:B118 description cell).

$ws.Range("A119").Value = 'let li = e.currentTarget.parentElement;
let id = li.getAttribute(''data-id'');'
$ws.Range("B119").Value = 'Изтриване на елемент от сървъра .'

$ws.Range("A120").Value = 'let id = li.getAttribute(''data-id'');'
$ws.Range("B120").Value = 'Взимане на атрибут на елемент'

# --- Formatting --------------------------------------------------------------
# B117/B118 hold the (merged) Bulgarian description -> centered + wrapped.
$ws.Range("B117:B118").HorizontalAlignment = -4108
$ws.Range("B117:B118").WrapText = $true

# A118/A119 hold wrapped code/description text (left aligned, wrap only).
$ws.Range("A118").WrapText = $true
$ws.Range("A119").WrapText = $true

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(117).RowHeight = 28.5
$ws.Rows.Item(118).RowHeight = 63
$ws.Rows.Item(119).RowHeight = 31.5

# --- Merge the description cell spanning the two code rows -------------------
$ws.Range("B117:B118").Merge() | Out-Null

# --- Selection / view state matches the saved workbook -----------------------
$ws.Range("B120").Select() | Out-Null
